$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7160.25
$ws.Range("I51").Value = 3996
$ws.Range("J51").Value = 8215
$ws.Range("K51").Value = 3996
$ws.Range("L51").Value = 8215
$ws.Range("M51").Value = -3512
$ws.Range("N51").Value = -9183
$ws.Range("H113").Value = 3806
$ws.Range("I113").Value = 3434
$ws.Range("J113").Value = 4085
$ws.Range("K113").Value = 3434
$ws.Range("L113").Value = 4085
$ws.Range("M113").Value = -180
$ws.Range("N113").Value = -10593
$ws.Range("H116").Value = 5350.143
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H137").Value = 3437.9285
$ws.Range("I137").Value = 4160.7837
$ws.Range("K137").Value = 12482.3511
$ws.Range("M137").Value = -9932.3511
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3041.9124
$ws.Range("I61").Value = 3041.9124
$ws.Range("K61").Value = 3041.9124
$ws.Range("M61").Value = -2829.9124
$ws.Range("H74").Value = 4615.9575
$ws.Range("I74").Value = 4875.7207
$ws.Range("J74").Value = 1823.5
$ws.Range("K74").Value = 4875.7207
$ws.Range("L74").Value = 1823.5
$ws.Range("M74").Value = -4001.7207
$ws.Range("N74").Value = -3571.5
$ws.Range("H77").Value = 4615.9575
$ws.Range("I77").Value = 4875.7207
$ws.Range("J77").Value = 1823.5
$ws.Range("K77").Value = 24378.6035
$ws.Range("L77").Value = 9117.5
$ws.Range("M77").Value = -20010.6035
$ws.Range("N77").Value = -17853.5
$ws.Range("H102").Value = 4089.0908
$ws.Range("I102").Value = 4484.4287
$ws.Range("K102").Value = 4484.4287
$ws.Range("M102").Value = -2862.4287
$ws.Range("H136").Value = 3041.9124
$ws.Range("I136").Value = 3041.9124
$ws.Range("K136").Value = 9125.7372
$ws.Range("M136").Value = -6575.7372
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 670.6842
$ws.Range("I80").Value = 490
$ws.Range("K80").Value = 490
$ws.Range("M80").Value = 508
$ws.Range("H83").Value = 670.6842
$ws.Range("I83").Value = 490
$ws.Range("K83").Value = 2450
$ws.Range("M83").Value = 2542
$ws.Range("H99").Value = 1597.75
$ws.Range("I99").Value = 1608.4
$ws.Range("J99").Value = 1544.5
$ws.Range("K99").Value = 1608.4
$ws.Range("L99").Value = 1544.5
$ws.Range("M99").Value = -110.4000000000001
$ws.Range("N99").Value = -4540.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2015.6875
$ws.Range("I16").Value = 1277.8
$ws.Range("K16").Value = 1277.8
$ws.Range("M16").Value = -990.8
$ws.Range("H86").Value = 11122444
$ws.Range("I86").Value = 16681416
$ws.Range("K86").Value = 16681416
$ws.Range("M86").Value = -16680293
$ws.Range("H89").Value = 11122444
$ws.Range("I89").Value = 16681416
$ws.Range("K89").Value = 83407080
$ws.Range("M89").Value = -83401464
$ws.Range("H99").Value = 8785.975
$ws.Range("J99").Value = 8809.071
$ws.Range("L99").Value = 8809.071
$ws.Range("N99").Value = -11805.071
$ws.Range("H105").Value = 3261.35
$ws.Range("I105").Value = 2078.5
$ws.Range("J105").Value = 4049.9167
$ws.Range("K105").Value = 2078.5
$ws.Range("L105").Value = 4049.9167
$ws.Range("M105").Value = -331.5
$ws.Range("N105").Value = -7543.9167
$ws.Range("H107").Value = 2000.2174
$ws.Range("I107").Value = 1694.85
$ws.Range("K107").Value = 1694.85
$ws.Range("M107").Value = 225.1500000000001
$ws.Range("H113").Value = 2015.6875
$ws.Range("I113").Value = 1277.8
$ws.Range("K113").Value = 1277.8
$ws.Range("M113").Value = 892.2
$ws.Range("H126").Value = 8785.975
$ws.Range("J126").Value = 8809.071
$ws.Range("L126").Value = 26427.213
$ws.Range("N126").Value = -31367.213
$ws.Range("H132").Value = 3565.5908
$ws.Range("J132").Value = 4344
$ws.Range("L132").Value = 13032
$ws.Range("N132").Value = -18092
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 9999.25
$ws.Range("I74").Value = 4999
$ws.Range("K74").Value = 14997
$ws.Range("M74").Value = -13936
$ws.Range("H77").Value = 9999.25
$ws.Range("I77").Value = 4999
$ws.Range("K77").Value = 44991
$ws.Range("M77").Value = -39687
$ws.Range("H109").Value = 2000
$ws.Range("I109").Value = 2000
$ws.Range("K109").Value = 6000
$ws.Range("M109").Value = -4960
$ws.Range("H113").Value = 1178.95
$ws.Range("I113").Value = 882.44446
$ws.Range("J113").Value = 1421.5454
$ws.Range("K113").Value = 2647.33338
$ws.Range("L113").Value = 4264.6362
$ws.Range("M113").Value = -477.33338
$ws.Range("N113").Value = -8604.636200000001
$ws.Range("H140").Value = 1724.8422
$ws.Range("I140").Value = 1724.8422
$ws.Range("K140").Value = 5174.5266
$ws.Range("M140").Value = 5.473399999999856
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1070
$ws.Range("I97").Value = 998.38464
$ws.Range("K97").Value = 998.38464
$ws.Range("M97").Value = -502.38464
$ws.Range("H122").Value = 4551.0713
$ws.Range("I122").Value = 3889.32
$ws.Range("K122").Value = 11667.96
$ws.Range("M122").Value = -9217.960000000001
$ws.Range("H126").Value = 5785.95
$ws.Range("J126").Value = 7007.75
$ws.Range("L126").Value = 21023.25
$ws.Range("N126").Value = -25963.25
$ws.Range("H136").Value = 62779.223
$ws.Range("J136").Value = 62779.223
$ws.Range("L136").Value = 188337.669
$ws.Range("N136").Value = -193437.669
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3927.4358
$ws.Range("J40").Value = 3440.9285
$ws.Range("L40").Value = 3440.9285
$ws.Range("N40").Value = -3712.9285
$ws.Range("H46").Value = 5853.727
$ws.Range("I46").Value = 5600
$ws.Range("K46").Value = 5600
$ws.Range("M46").Value = -5412
$ws.Range("H93").Value = 22223534
$ws.Range("I93").Value = 28572542
$ws.Range("J93").Value = 2002
$ws.Range("K93").Value = 28572542
$ws.Range("L93").Value = 2002
$ws.Range("M93").Value = -28571294
$ws.Range("N93").Value = -4498
$ws.Range("H100").Value = 90911064
$ws.Range("J100").Value = 2343.5715
$ws.Range("L100").Value = 2343.5715
$ws.Range("N100").Value = -3425.5715
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 28874.5
$ws.Range("J70").Value = 28874.5
$ws.Range("L70").Value = 28874.5
$ws.Range("N70").Value = -29504.5
$ws.Range("H73").Value = 28874.5
$ws.Range("J73").Value = 28874.5
$ws.Range("L73").Value = 28874.5
$ws.Range("N73").Value = -31058.5
